$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) cells so numeric-looking strings (e.g. trailing
# zeros, multi-dot thousands separators) are preserved verbatim instead of being
# auto-coerced into numbers by Excel.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D9", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D29", "D30", "D31", "D33", "D35", "D36", "D37", "D38", "D41", "D42", "D44", "D45", "D47", "D48", "D50", "D51")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '58.327.57'
$ws.Range('E2').Value = '  +3.30%  '
$ws.Range('D3').Value = '2.365.08'
$ws.Range('E3').Value = '  +1.56%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').Value = '543.55'
$ws.Range('E5').Value = '  +6.04%  '
$ws.Range('D6').Value = '135.00'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  +0.88%  '
$ws.Range('D9').Value = '2.363.17'
$ws.Range('E9').Value = '  +1.29%  '
$ws.Range('E10').Value = '  +1.90%  '
$ws.Range('E11').Value = '  +1.22%  '
$ws.Range('E12').Value = '  +2.22%  '
$ws.Range('E13').Value = '  +5.22%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = '23.66'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.752.51'
$ws.Range('E15').Value = '  +0.32%  '
$ws.Range('D16').Value = '58.198.58'
$ws.Range('E16').Value = '  +3.09%  '
$ws.Range('E17').Value = '  +1.17%  '
$ws.Range('D18').Value = '2.357.69'
$ws.Range('E18').Value = '  +1.45%  '
$ws.Range('D19').Value = '10.59'
$ws.Range('E19').Value = '  +1.42%  '
$ws.Range('D20').Value = '334.81'
$ws.Range('E20').Value = '  +2.49%  '
$ws.Range('D21').Value = '4.22'
$ws.Range('E21').Value = '  +1.99%  '
$ws.Range('D22').Value = '6.75'
$ws.Range('E22').Value = '  +0.42%  '
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('D24').Value = '62.09'
$ws.Range('E24').Value = '  +1.02%  '
$ws.Range('E25').Value = '  +4.84%  '
$ws.Range('D26').Value = '8.51'
$ws.Range('E26').Value = '  -3.41%  '
$ws.Range('D27').Value = '0.994'
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('E28').Value = '  +10.72%  '
$ws.Range('D29').Value = '1.76'
$ws.Range('E29').Value = '  +5.20%  '
$ws.Range('D30').Value = '170.67'
$ws.Range('E30').Value = '  +1.86%  '
$ws.Range('D31').Value = '0.0₃0740'
$ws.Range('E31').Value = '  +2.96%  '
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('D33').Value = '18.59'
$ws.Range('E33').Value = '  +1.26%  '
$ws.Range('E34').Value = '  +14.76%  '
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').Value = '4.18'
$ws.Range('E37').Value = '  +6.55%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '1.27'
$ws.Range('E38').Value = '  +0.57%  '
$ws.Range('E39').Value = '  +5.82%  '
$ws.Range('E40').Value = '  +2.32%  '
$ws.Range('D41').Value = '150.00'
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('D42').Value = '0.381'
$ws.Range('E42').Value = '  +1.98%  '
$ws.Range('E43').Value = '  +2.07%  '
$ws.Range('D44').Value = '285.44'
$ws.Range('E44').Value = '  +2.69%  '
$ws.Range('D45').Value = '19.38'
$ws.Range('E45').Value = '  +6.67%  '
$ws.Range('E46').Value = '  +0.39%  '
$ws.Range('D47').Value = '0.0508'
$ws.Range('E47').Value = '  +2.90%  '
$ws.Range('D48').Value = '0.564'
$ws.Range('E48').Value = '  +1.46%  '
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('D50').Value = '17.71'
$ws.Range('E50').Value = '  +3.53%  '
$ws.Range('D51').Value = '0.382'
$ws.Range('E51').Value = '  +0.80%  '

# Restore default (Normal) style on the Price cells now that the text value is set,
# so no stray style index lingers on cells that had none originally.
foreach ($c in $priceCells) {
    $ws.Range($c).Style = "Normal"
}
